# Calibration of energy use modeling by renovation level:
# Divide all TH values in column B (rows 2-452) by 3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 452
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $current = $cell.Value()
    $cell.Value = $current / 3
}
